# Reproduces the commit:
#   Revert "Revert "added requirement openpyxl""
#
# Net effect vs. before.xlsx:
#   1. sheet "test_file" (Sheet1): the "sex" column (E) is removed entirely,
#      so the former "customer_type" column (F) shifts left into column E.
#   2. sheet "Sheet2": four more sample rows (4-7) are appended, repeating
#      the existing 1,2,3 / 2,3,4 pattern.
#   3. The active/selected sheet goes back to "test_file" (first sheet),
#      with the cursor on I8; Sheet2's own remembered selection becomes F16.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Drop the "sex" column (column E) from the first sheet.
$ws1.Columns.Item(5).Delete()

# 2) Append the extra rows to Sheet2.
$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 3

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = 3
$ws2.Range("C5").Value = 4

$ws2.Range("A6").Value = 1
$ws2.Range("B6").Value = 2
$ws2.Range("C6").Value = 3

$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = 3
$ws2.Range("C7").Value = 4

# 3) Selection / active-sheet bookkeeping: Sheet2 keeps a remembered
#    selection of F16, but the workbook re-opens on the first sheet with
#    I8 selected.
$ws2.Range("F16").Select()

$ws1.Activate()
$ws1.Range("I8").Select()
